# Apply the vocabulary.xlsx update:
#  - insert a new "dct:creator" row (additional ORCID) right after the
#    existing dct:creator row, pushing the metadata/definition/term rows
#    down by one
#  - rename the two sample terms and append the real NICEST-2 subject terms
#    generated from the Google Sheet / sheet2rdf workflow

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 11 (shifts old rows 11-18 down to 12-19)
$ws.Rows(11).Insert()

# New row 11: a second dct:creator entry (ORCID), no note in column C
$ws.Cells.Item(11, 1).Value = "dct:creator"
$ws.Cells.Item(11, 2).Value = "https://orcid.org/0000-0001-8497-1661"

# Rename the two placeholder sample terms (now at rows 18 & 19 after the insert)
$ws.Cells.Item(18, 2).Value = "oceanography"
$ws.Cells.Item(19, 2).Value = "numerical modelling"

# Append the remaining NICEST-2 subject terms generated from the sheet
$ws.Cells.Item(20, 1).Value = "nicest-2-subjects:10002"
$ws.Cells.Item(20, 2).Value = "forecast"

$ws.Cells.Item(21, 1).Value = "nicest-2-subjects:10003"
$ws.Cells.Item(21, 2).Value = "observations"

$ws.Cells.Item(22, 1).Value = "nicest-2-subjects:10004"
$ws.Cells.Item(22, 2).Value = "meteorology"

$ws.Cells.Item(23, 1).Value = "nicest-2-subjects:10005"
$ws.Cells.Item(23, 2).Value = "timeseries"

$ws.Cells.Item(24, 1).Value = "nicest-2-subjects:10006"
$ws.Cells.Item(24, 2).Value = "hydrology"

$ws.Cells.Item(25, 1).Value = "nicest-2-subjects:10007"
$ws.Cells.Item(25, 2).Value = "atmosphere"

$ws.Cells.Item(26, 1).Value = "nicest-2-subjects:10008"
$ws.Cells.Item(26, 2).Value = "sea ice"
$ws.Cells.Item(26, 6).Value = "nicest-2-subjects:10000"
$ws.Cells.Item(26, 9).Value = "to be reviewed"

$ws.Cells.Item(27, 1).Value = "nicest-2-subjects:10009"
$ws.Cells.Item(27, 2).Value = "cryosphere"

$ws.Cells.Item(28, 1).Value = "nicest-2-subjects:10010"

$ws.Cells.Item(29, 1).Value = "nicest-2-subjects:10011"

$ws.Cells.Item(30, 1).Value = "nicest-2-subjects:10012"
